# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which carry the same data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 1167
    3  = 108
    4  = 1598
    5  = 613
    8  = 11416
    9  = 23
    11 = 447
    12 = 351
    14 = 791
    15 = 12356
    16 = 13026
    20 = 35
    21 = 13
    24 = 99
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
